$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 1.73
$ws.Range("O2").Value = 2.1
$ws.Range("G4").Value = 1.7
$ws.Range("I4").Value = 5.5
$ws.Range("G7").Value = 1.9
$ws.Range("I7").Value = 3.7
$ws.Range("V7").Value = 8.5
$ws.Range("AA7").Value = 6.5
$ws.Range("AB7").Value = 13
$ws.Range("AD7").Value = 12
$ws.Range("AE7").Value = 21
$ws.Range("G8").Value = 2.88
$ws.Range("I8").Value = 2.3
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 11
$ws.Range("AG8").Value = 21
$ws.Range("G9").Value = 1.7
$ws.Range("H9").Value = 3.7
$ws.Range("I9").Value = 5
$ws.Range("N9").Value = 1.6
$ws.Range("O9").Value = 2.3
$ws.Range("T9").Value = 9.5
$ws.Range("AB9").Value = 12
$ws.Range("G13").Value = 2.55
$ws.Range("I13").Value = 2.75
$ws.Range("K13").Value = 8.5
$ws.Range("N13").Value = 2.1
$ws.Range("O13").Value = 1.7
$ws.Range("U13").Value = 12
$ws.Range("V13").Value = 10
$ws.Range("W13").Value = 26
$ws.Range("X13").Value = 21
$ws.Range("AE13").Value = 13
$ws.Range("AF13").Value = 11
$ws.Range("AG13").Value = 29
$ws.Range("AH13").Value = 23
$ws.Range("J15").Value = 1.1
$ws.Range("K15").Value = 7
$ws.Range("N15").Value = 2.6
$ws.Range("O15").Value = 1.48
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 10
$ws.Range("X17").Value = 21
$ws.Range("AB17").Value = 13
$ws.Range("G20").Value = 2.1
$ws.Range("I20").Value = 3.8
$ws.Range("J20").Value = 1.08
$ws.Range("K20").Value = 8
$ws.Range("U20").Value = 9.5
$ws.Range("V20").Value = 9.5
$ws.Range("W20").Value = 19
$ws.Range("X20").Value = 19
$ws.Range("AD20").Value = 10
$ws.Range("AF20").Value = 13
$ws.Range("J21").Value = 1.06
$ws.Range("K21").Value = 10
$ws.Range("L21").Value = 1.25
$ws.Range("M21").Value = 3.75
$ws.Range("N21").Value = 1.83
$ws.Range("O21").Value = 2.03
$ws.Range("Y21").Value = 29
$ws.Range("Z21").Value = 10
$ws.Range("AB21").Value = 19
$ws.Range("G24").Value = 1.7
$ws.Range("H24").Value = 3.75
$ws.Range("J24").Value = 1.04
$ws.Range("K24").Value = 13
$ws.Range("N24").Value = 1.83
$ws.Range("O24").Value = 2.03
$ws.Range("P24").Value = 1.36
$ws.Range("Q24").Value = 3
$ws.Range("T24").Value = 7.5
$ws.Range("X24").Value = 13
$ws.Range("G25").Value = 1.53
$ws.Range("H25").Value = 3.8
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 1.07
$ws.Range("K25").Value = 8.5
$ws.Range("L25").Value = 1.33
$ws.Range("M25").Value = 3.25
$ws.Range("R25").Value = 2.1
$ws.Range("S25").Value = 1.67
$ws.Range("T25").Value = 6
$ws.Range("U25").Value = 6.5
$ws.Range("AC25").Value = 67
$ws.Range("G26").Value = 2.3
$ws.Range("J26").Value = 1.11
$ws.Range("K26").Value = 6.5
$ws.Range("Z26").Value = 6.5
$ws.Range("AB26").Value = 17
$ws.Range("J27").Value = 1.11
$ws.Range("K27").Value = 6.5
$ws.Range("G29").Value = 3.8
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 1.07
$ws.Range("K29").Value = 9
$ws.Range("L29").Value = 1.36
$ws.Range("M29").Value = 3
$ws.Range("N29").Value = 2.2
$ws.Range("O29").Value = 1.65
$ws.Range("V29").Value = 13
$ws.Range("AC29").Value = 51
$ws.Range("AD29").Value = 6.5
$ws.Range("AH29").Value = 19
$ws.Range("G30").Value = 21
$ws.Range("H30").Value = 6.7
$ws.Range("I30").Value = 1.09
$ws.Range("N30").Value = 1.32
$ws.Range("O30").Value = 3.1
$ws.Range("R30").Value = 2.24
$ws.Range("S30").Value = 1.58
$ws.Range("T30").Value = 60
$ws.Range("U30").Value = 250
$ws.Range("V30").Value = 60
$ws.Range("Y30").Value = 150
$ws.Range("Z30").Value = 19.5
$ws.Range("AA30").Value = 14.5
$ws.Range("AB30").Value = 29
$ws.Range("AC30").Value = 110
$ws.Range("AD30").Value = 8.5
$ws.Range("AE30").Value = 5.8
$ws.Range("AF30").Value = 9.75
$ws.Range("AG30").Value = 5.5
$ws.Range("AH30").Value = 9
$ws.Range("AI30").Value = 27
$ws.Range("G36").Value = 1.75
$ws.Range("H36").Value = 3.75
$ws.Range("I36").Value = 4.2
$ws.Range("K36").Value = 12
$ws.Range("L36").Value = 1.25
$ws.Range("M36").Value = 3.75
$ws.Range("N36").Value = 1.83
$ws.Range("O36").Value = 2.03
$ws.Range("P36").Value = 1.36
$ws.Range("Q36").Value = 3
$ws.Range("R36").Value = 1.73
$ws.Range("S36").Value = 2
$ws.Range("T36").Value = 7.5
$ws.Range("U36").Value = 9
$ws.Range("W36").Value = 15
$ws.Range("X36").Value = 15
$ws.Range("Z36").Value = 12
$ws.Range("AA36").Value = 7.5
$ws.Range("AB36").Value = 15
$ws.Range("AE36").Value = 21
$ws.Range("AG36").Value = 41
$ws.Range("N38").Value = 1.85
$ws.Range("O38").Value = 2
$ws.Range("N41").Value = 1.88
$ws.Range("O41").Value = 1.98
$ws.Range("G46").Value = 2.35
$ws.Range("I46").Value = 2.75
$ws.Range("J46").Value = 1.05
$ws.Range("K46").Value = 8.5
$ws.Range("N46").Value = 1.83
$ws.Range("O46").Value = 1.98
$ws.Range("Y46").Value = 26
$ws.Range("AD46").Value = 10
$ws.Range("AG46").Value = 29
$ws.Range("G47").Value = 2.7
$ws.Range("H47").Value = 3.3
$ws.Range("I47").Value = 2.45
$ws.Range("J47").Value = 1.05
$ws.Range("K47").Value = 8.5
$ws.Range("T47").Value = 9
$ws.Range("U47").Value = 13
$ws.Range("W47").Value = 26
$ws.Range("X47").Value = 21
$ws.Range("Y47").Value = 29
$ws.Range("AD47").Value = 8.5
$ws.Range("AE47").Value = 12
$ws.Range("AG47").Value = 23
$ws.Range("AH47").Value = 21
$ws.Range("G48").Value = 1.91
$ws.Range("H48").Value = 4.1
$ws.Range("I48").Value = 3.2
$ws.Range("J48").Value = 23
$ws.Range("K48").Value = 1.03
$ws.Range("L48").Value = 1.11
$ws.Range("M48").Value = 6
$ws.Range("R48").Value = 1.4
$ws.Range("S48").Value = 2.75
$ws.Range("U48").Value = 13
$ws.Range("V48").Value = 9.5
$ws.Range("W48").Value = 19
$ws.Range("AC48").Value = 29
$ws.Range("AD48").Value = 17
$ws.Range("AE48").Value = 21
$ws.Range("AF48").Value = 12
$ws.Range("AG48").Value = 34
$ws.Range("AH48").Value = 21
$ws.Range("G49").Value = 1.95
$ws.Range("H49").Value = 3.6
$ws.Range("I49").Value = 3.7
$ws.Range("P49").Value = 1.36
$ws.Range("Q49").Value = 3
$ws.Range("U49").Value = 10
$ws.Range("W49").Value = 17
$ws.Range("X49").Value = 15
$ws.Range("Z49").Value = 12
$ws.Range("G53").Value = 2.15
$ws.Range("H53").Value = 2.92
$ws.Range("I53").Value = 3.7
$ws.Range("J53").Value = 1.07
$ws.Range("Q53").Value = 2.5
$ws.Range("S53").Value = 2.1
$ws.Range("T53").Value = 7.6
$ws.Range("U53").Value = 11
$ws.Range("V53").Value = 8.25
$ws.Range("W53").Value = 22
$ws.Range("X53").Value = 17
$ws.Range("AA53").Value = 5.7
$ws.Range("AB53").Value = 12
$ws.Range("AD53").Value = 10.75
$ws.Range("AF53").Value = 11.75
$ws.Range("AI53").Value = 35
$ws.Range("AJ53").Value = 400
$ws.Range("G54").Value = 7.5
$ws.Range("I54").Value = 1.38
$ws.Range("R54").Value = 1.67
$ws.Range("S54").Value = 2.1
$ws.Range("T54").Value = 26
$ws.Range("AF54").Value = 8.5
$ws.Range("G56").Value = 2.3
$ws.Range("H56").Value = 2.75
$ws.Range("I56").Value = 3.45
$ws.Range("L56").Value = 1.42
$ws.Range("M56").Value = 2.47
$ws.Range("N56").Value = 2.2
$ws.Range("O56").Value = 1.52
$ws.Range("P56").Value = 1.55
$ws.Range("Q56").Value = 2.15
$ws.Range("R56").Value = 1.83
$ws.Range("S56").Value = 1.78
$ws.Range("T56").Value = 6.3
$ws.Range("U56").Value = 10.5
$ws.Range("V56").Value = 9
$ws.Range("W56").Value = 25
$ws.Range("X56").Value = 22
$ws.Range("Y56").Value = 35
$ws.Range("Z56").Value = 6.7
$ws.Range("AA56").Value = 5.4
$ws.Range("AB56").Value = 14.5
$ws.Range("AC56").Value = 80
$ws.Range("AD56").Value = 8.75
$ws.Range("AE56").Value = 18
$ws.Range("AF56").Value = 11.75
$ws.Range("AG56").Value = 55
$ws.Range("AH56").Value = 35
$ws.Range("AI56").Value = 45
$ws.Range("AJ56").Value = 700
